$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H16").Value = 14331.667
$ws.Range("I16").Value = 7992.5
$ws.Range("J16").Value = 27010
$ws.Range("K16").Value = 7992.5
$ws.Range("L16").Value = 27010
$ws.Range("M16").Value = -7762.5
$ws.Range("N16").Value = -27470
$ws.Range("H39").Value = 434.46155
$ws.Range("I39").Value = 261
$ws.Range("J39").Value = 824.75
$ws.Range("K39").Value = 783
$ws.Range("L39").Value = 2474.25
$ws.Range("M39").Value = -487
$ws.Range("N39").Value = -3066.25
$ws.Range("H41").Value = 250
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H43").Value = 2750
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2516
$ws.Range("N51").ClearContents()
$ws.Range("H53").Value = 951
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H62").Value = 7053.3335
$ws.Range("I62").Value = 6413.5
$ws.Range("K62").Value = 6413.5
$ws.Range("M62").Value = -5789.5
$ws.Range("H65").Value = 7053.3335
$ws.Range("I65").Value = 6413.5
$ws.Range("K65").Value = 32067.5
$ws.Range("M65").Value = -28947.5
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1685
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -908
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 15000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -26232
$ws.Range("H106").Value = 1700
$ws.Range("I106").Value = 1700
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1700
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1069
$ws.Range("N106").ClearContents()
$ws.Range("H135").Value = 882.75
$ws.Range("I135").Value = 882.75
$ws.Range("K135").Value = 7944.75
$ws.Range("M135").Value = -5409.75

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 1468.1538
$ws.Range("I32").Value = 1257.1666
$ws.Range("K32").Value = 1257.1666
$ws.Range("M32").Value = -970.1666
$ws.Range("H44").Value = 19500
$ws.Range("J44").Value = 19500
$ws.Range("L44").Value = 19500
$ws.Range("N44").Value = -20476
$ws.Range("H55").Value = 20000
$ws.Range("J55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("N55").Value = -20630

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1069.9166
$ws.Range("I20").Value = 878.8570999999999
$ws.Range("K20").Value = 878.8570999999999
$ws.Range("M20").Value = -631.8570999999999
$ws.Range("H37").Value = 3153.5
$ws.Range("I37").Value = 378.4
$ws.Range("J37").Value = 17029
$ws.Range("K37").Value = 378.4
$ws.Range("L37").Value = 17029
$ws.Range("M37").Value = -241.4
$ws.Range("N37").Value = -17303
$ws.Range("H80").Value = 236.3
$ws.Range("I80").Value = 85.666664
$ws.Range("K80").Value = 85.666664
$ws.Range("M80").Value = 912.333336
$ws.Range("H83").Value = 236.3
$ws.Range("I83").Value = 85.666664
$ws.Range("K83").Value = 428.33332
$ws.Range("M83").Value = 4563.66668
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 45000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 3355.6667
$ws.Range("I58").Value = 1965
$ws.Range("J58").Value = 4349
$ws.Range("K58").Value = 1965
$ws.Range("L58").Value = 4349
$ws.Range("M58").Value = -1762
$ws.Range("N58").Value = -4755
$ws.Range("H134").Value = 1570.6666
$ws.Range("I134").Value = 1570.6666
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4711.9998
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2176.9998
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3355.6667
$ws.Range("I136").Value = 1965
$ws.Range("J136").Value = 4349
$ws.Range("K136").Value = 5895
$ws.Range("L136").Value = 13047
$ws.Range("M136").Value = -3345
$ws.Range("N136").Value = -18147

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 387.58334
$ws.Range("I2").Value = 188.25
$ws.Range("J2").Value = 786.25
$ws.Range("K2").Value = 1129.5
$ws.Range("L2").Value = 4717.5
$ws.Range("M2").Value = -1016.5
$ws.Range("N2").Value = -4943.5
$ws.Range("H34").Value = 1169
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1169
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3507
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3675
$ws.Range("H39").Value = 4088.4285
$ws.Range("J39").Value = 4249.077
$ws.Range("L39").Value = 12747.231
$ws.Range("N39").Value = -13335.231
$ws.Range("H129").Value = 1015
$ws.Range("J129").Value = 1000
$ws.Range("L129").Value = 3000
$ws.Range("N129").Value = -13000
$ws.Range("H131").Value = 3333.3333

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H70").Value = 6076
$ws.Range("I70").Value = 7258
$ws.Range("K70").Value = 7258
$ws.Range("M70").Value = -6988
$ws.Range("H73").Value = 6076
$ws.Range("I73").Value = 7258
$ws.Range("K73").Value = 7258
$ws.Range("M73").Value = -6322

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 9999
$ws.Range("I136").Value = 9999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 29997
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -27447
$ws.Range("N136").ClearContents()

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4490
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
